$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCpUC")
$ws.Rows.Item(2).Delete()
[void]$ws.Range("A2:XFD2").Select()
